$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to text format so that values such as
# "9.00", "3.50", "0.0000271" and thousand-separated prices like "63.909.92"
# are written verbatim instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Simple price / 1h-volume updates -------------------------------------
$ws.Range("D2").Value = '63.909.92'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '3.337.18'
$ws.Range("E3").Value = '  -3.80%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '556.15'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").Value = '174.87'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  -2.36%  '
$ws.Range("D8").Value = '3.331.25'
$ws.Range("E8").Value = '  -3.85%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '0.622'
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("D11").Value = '0.161'
$ws.Range("E11").Value = '  +5.52%  '
$ws.Range("D12").Value = '53.94'
$ws.Range("E12").Value = '  +1.70%  '
$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").Value = '9.00'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '3.880.00'
$ws.Range("E15").Value = '  -3.51%  '
$ws.Range("E16").Value = '  +1.23%  '
$ws.Range("E17").Value = '  -2.25%  '
$ws.Range("D18").Value = '3.347.34'
$ws.Range("E18").Value = '  -3.38%  '
$ws.Range("D19").Value = '11.84'
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").Value = '63.900.19'
$ws.Range("E20").Value = '  -2.37%  '
$ws.Range("D21").Value = '0.978'
$ws.Range("E21").Value = '  -0.71%  '
$ws.Range("D22").Value = '433.85'
$ws.Range("E22").Value = '  +5.83%  '
$ws.Range("D23").Value = '4.56'
$ws.Range("E23").Value = '  +10.02%  '
$ws.Range("D24").Value = '4.11'
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").Value = '84.26'
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").Value = '13.16'
$ws.Range("E26").Value = '  +4.39%  '
$ws.Range("D29").Value = '8.74'
$ws.Range("E29").Value = '  -1.26%  '
$ws.Range("D30").Value = '29.47'
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").Value = '6.53'
$ws.Range("E31").Value = '  +4.81%  '
$ws.Range("D32").Value = '591.90'
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("D33").Value = '11.46'
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("E34").Value = '  -0.88%  '
$ws.Range("D35").Value = '58.61'
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").Value = '0.144'
$ws.Range("E37").Value = '  -3.39%  '
$ws.Range("D38").Value = '3.50'
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("D41").Value = '0.365'
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("D42").Value = '3.112.03'
$ws.Range("E42").Value = '  -7.21%  '
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = '2.84'
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("D47").Value = '2.43'
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("E49").Value = '  -3.01%  '
$ws.Range("E50").Value = '  -1.57%  '
$ws.Range("D51").Value = '133.10'
$ws.Range("E51").Value = '  -3.12%  '

# --- Ranking swaps: RenderToken/ImmutableX (rows 27-28) -------------------
# --- and InjectiveProtocol/PEPE (rows 39-40) -------------------------------
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '10.63'
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = '2.83'
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").Value = '35.33'
$ws.Range("E39").Value = '  -3.02%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0749'
$ws.Range("E40").Value = '  -4.11%  '

# Restore the default (unstyled) look for the Price/Volume columns so the
# text-format coercion above does not leave a stray cell style behind.
$ws.Range("D2:E51").Style = "Normal"
